# Auto-generated Excel COM-interop script to update cryptos.xlsx data
# Applies per-cell text/value updates to match the target snapshot of crypto
# prices / volume percentages, plus a few row re-orderings (coin swaps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cell updates (strings that Excel will not reinterpret as numbers) ---
$textUpdates = @{
    'D2' = '60.119.78'
    'E2' = '  +4.42%  '
    'D3' = '2.587.60'
    'E3' = '  +5.58%  '
    'E4' = '  +0.17%  '
    'E5' = '  +3.11%  '
    'E6' = '  -0.01%  '
    'E7' = '  +0.11%  '
    'E8' = '  -4.66%  '
    'D9' = '2.622.05'
    'E9' = '  +5.78%  '
    'E10' = '  +3.68%  '
    'E11' = '  +2.76%  '
    'E12' = '  +2.16%  '
    'E13' = '  +1.13%  '
    'D14' = '3.068.92'
    'E14' = '  +6.53%  '
    'D15' = '60.248.61'
    'E15' = '  +4.72%  '
    'E16' = '  +4.22%  '
    'E17' = '  +3.76%  '
    'D18' = '2.615.63'
    'E18' = '  +5.52%  '
    'E19' = '  +2.74%  '
    'E20' = '  +5.31%  '
    'E21' = '  +2.88%  '
    'E22' = '  +2.69%  '
    'E23' = '  +0.25%  '
    'B24' = 'Litecoin'
    'C24' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'E24' = '  +2.38%  '
    'B25' = 'Polygon'
    'C25' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E25' = '  +4.35%  '
    'B26' = 'WrappedeETH'
    'C26' = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    'D26' = '2.720.51'
    'E26' = '  +5.92%  '
    'E27' = '  +2.47%  '
    'E28' = '  +0.04%  '
    'E29' = '  +5.66%  '
    'E30' = '  +2.06%  '
    'E31' = '  +0.12%  '
    'E32' = '  +3.32%  '
    'E33' = '  +2.69%  '
    'E34' = '  +1.91%  '
    'E35' = '  +6.77%  '
    'E36' = '  +4.59%  '
    'E37' = '  +5.52%  '
    'E38' = '  +1.52%  '
    'E39' = '  +6.56%  '
    'B40' = 'SuiNetwork'
    'C40' = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
    'E40' = '  +27.75%  '
    'B41' = 'Stacks'
    'C41' = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    'E41' = '  +5.44%  '
    'B42' = 'Bittensor'
    'C42' = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    'E42' = '  +5.86%  '
    'E43' = '  +3.91%  '
    'E44' = '  +2.97%  '
    'E45' = '  +5.55%  '
    'E46' = '  -0.07%  '
    'E47' = '  -0.12%  '
    'E48' = '  +9.91%  '
    'E49' = '  +5.86%  '
    'D50' = '2.045.53'
    'E50' = '  +7.37%  '
    'E51' = '  +1.79%  '
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}

# --- Numeric-looking text updates (values such as '1.00' or '0.0566' that Excel
#     would otherwise auto-convert to a Number and strip trailing zeros). These
#     are forced to remain Text cells, matching the original inlineStr cells. ---
$numericLookingUpdates = @{
    'D4' = '1.00'
    'D5' = '506.24'
    'D6' = '155.40'
    'D7' = '0.997'
    'D8' = '0.584'
    'D10' = '6.45'
    'D12' = '0.341'
    'D16' = '21.66'
    'D20' = '342.41'
    'D21' = '10.37'
    'D22' = '6.11'
    'D23' = '0.998'
    'D24' = '59.87'
    'D25' = '0.421'
    'D28' = '0.993'
    'D30' = '7.49'
    'D32' = '156.03'
    'D33' = '19.30'
    'D34' = '1.56'
    'D35' = '5.69'
    'D36' = '3.97'
    'D40' = '0.838'
    'D41' = '1.47'
    'D42' = '300.38'
    'D43' = '35.57'
    'D44' = '0.624'
    'D45' = '0.0566'
    'D47' = '0.992'
    'D48' = '19.76'
    'D49' = '4.92'
}

foreach ($cell in $numericLookingUpdates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $numericLookingUpdates[$cell]
    $rng.Style = "Normal"
}

Write-Host "Updated $($textUpdates.Count + $numericLookingUpdates.Count) cells"